# Update hourly close/high/low rates for SUNTV sheet - closing dates were
# wrong for algo, so the High/Low/LTP/PREV figures for the affected rows
# need to be corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("SUNTV")

# Row 7 (daily summary row)
$ws.Range("F7").Value = 678.6
$ws.Range("G7").Value = 696
$ws.Range("H7").Value = 672.5
$ws.Range("I7").Value = 693.85
$ws.Range("J7").Value = 677.95

# Row 9
$ws.Range("G9").Value = 680.85
$ws.Range("H9").Value = 665.4
$ws.Range("I9").Value = 676.6

# Row 10
$ws.Range("G10").Value = 683.2
$ws.Range("H10").Value = 675.9
$ws.Range("I10").Value = 681.1

# Row 11
$ws.Range("G11").Value = 683.6
$ws.Range("H11").Value = 680.05
$ws.Range("I11").Value = 680.35

# Row 12
$ws.Range("G12").Value = 686.4
$ws.Range("H12").Value = 678
$ws.Range("I12").Value = 686.35

# Row 13
$ws.Range("G13").Value = 687.4
$ws.Range("H13").Value = 683.8
$ws.Range("I13").Value = 684.45

# Row 14
$ws.Range("G14").Value = 690.15
$ws.Range("H14").Value = 684
$ws.Range("I14").Value = 689.7

# Row 15
$ws.Range("G15").Value = 690
$ws.Range("H15").Value = 687.35
$ws.Range("I15").Value = 688.2

# Row 16
$ws.Range("G16").Value = 689.7
$ws.Range("H16").Value = 687.55
$ws.Range("I16").Value = 689.25

# Row 17
$ws.Range("G17").Value = 689.55
$ws.Range("H17").Value = 685.45
$ws.Range("I17").Value = 687.4

# Row 18
$ws.Range("G18").Value = 689.5
$ws.Range("H18").Value = 683.6
$ws.Range("I18").Value = 686.5

# Row 19
$ws.Range("G19").Value = 686.85
$ws.Range("H19").Value = 684.05
$ws.Range("I19").Value = 685.7

# Row 20
$ws.Range("G20").Value = 691.5
$ws.Range("H20").Value = 685.5
$ws.Range("I20").Value = 691.4

# Row 21
$ws.Range("G21").Value = 696
$ws.Range("H21").Value = 690.2
$ws.Range("I21").Value = 693

$wb.Save()
